$d = $word.ActiveDocument

# Locate the paragraph that ends "...adequately secure the file system." -
# the new content is appended to it, and several new paragraphs are
# inserted between it and the blank paragraph that follows (and precedes
# the section break).
$anchor = $d.Content
$anchor.Find.Execute("adequately secure the file system.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetPara = $anchor.Paragraphs(1).Range
$blankPara = $targetPara.Next(4, 1)

$r = $d.Range($targetPara.Start, $blankPara.End)

$xml = '<w:p w14:paraId="27BD9824" w14:textId="5326DB49" w:rsidR="00FC7DD9" w:rsidRPr="00FC7DD9" w:rsidRDefault="001D385F" w:rsidP="002D0078" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r w:rsidRPr="001D385F">
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">In contrast to these security mechanisms, Microsoft has recently added a number of additional capabilities to the Windows 10 OS to adequately secure the file system. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">They consist of Credential Guard, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Device Guard, and Secure Boot. Malicious </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>programmes</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> cannot load during the startup process thanks to a technology called Safe Boot. Administrators have the option to limit which </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>programmes</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> can be run on a device thanks to Device Guard.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Credential Guard, which isolates user credentials from hostile attackers, offers a last degree of safety. It secures and isolates the sensitive data it keeps, making it impossible for hackers to access them.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>To further improve the security of their file systems, new security measures have been added to subsequent versions of Windows 10. On Windows, these features consist of Secure Boot, Device Guard, and Credential Guard.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">By ensuring that only </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>authorised</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> kernels are booted during the boot process for Windows, Secure Boot prevents malicious malware from loading. Administrators can limit which applications can run on a device with Device Guard. By separating user credentials from hostile actors, Credential Guard adds an extra degree of security.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Device Guard, a feature supported by Windows 10, protects your device from harmful </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>programmes</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> by stopping it from running on the device. Credential Guard is an additional feature that aids in preventing device theft of sensitive credentials.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>A number of security improvements have been made to the most recent version of Windows 10 to help safeguard your device from harmful software. A more thorough threat protection system, including real-time protection, automated updates, enhanced security software, and a secure processor, has been added to Windows Security. User account control, BitLocker</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">firewall and network security, secure boot and Exploit Guard are </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>some</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>additional</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>we get in</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> Windows 10.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p w14:paraId="2F825580" w14:textId="77777777" w:rsidR="00FC7DD9" w:rsidRPr="00FC7DD9" w:rsidRDefault="00FC7DD9" w:rsidP="00FC7DD9">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'
$r.InsertXML($xml)
